$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) "Ativacao:" date string 01/01/2012 -> 01/01/2023 -------------------
# The original value is plain text (shared string), not a real date. Typing
# a dd/mm/yyyy-looking literal via .Value would get auto-converted by Excel
# into a date serial number, so instead we write it as a formula result and
# then flatten the formula to its static text value via copy/paste-values.
# This keeps the cell as a shared text string (t="s") and preserves the
# pre-existing cell style (the B column "normal" style, the C column "red"
# style used for the values highlighted as changed).
$ws.Range("B8").Formula = '="01/01/2023"'
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("C8").Formula = '="01/01/2023"'
$ws.Range("C8").Copy()
$ws.Range("C8").PasteSpecial(-4163)   # xlPasteValues

# Row 15 ("Programa:") happens to reuse the very same shared string as the
# "Ativacao:" date (pre-existing quirk of this sheet), so it must show the
# same updated text too.
$ws.Range("B15").Formula = '="01/01/2023"'
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("C15").Formula = '="01/01/2023"'
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)  # xlPasteValues

# --- 2) New English "Objectives" paragraph on row 11 -----------------------
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats (normal/black style)
$ws.Range("B11").Value = "Rheology is the science that studies the flow of materials. Your knowledge is necessary to understand the processes of forming materials. The course aims to provide students with the basic and applied concepts of rheology and familiarize them with experimental methods for evaluating the rheological properties of materials."

$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)  # xlPasteFormats (red style)
$ws.Range("C11").Value = "Rheology is the science that studies the flow of materials. Your knowledge is necessary to understand the processes of forming materials. The course aims to provide students with the basic and applied concepts of rheology and familiarize them with experimental methods for evaluating the rheological properties of materials."

# --- 3) New English "Short syllabus" paragraph on row 14 --------------------
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B14").Value = "Flow of Newtonian and non-Newtonian fluids. Viscosity and rheometry. viscoelasticity. Applications."

$ws.Range("C10").Copy()
$ws.Range("C14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C14").Value = "Flow of Newtonian and non-Newtonian fluids. Viscosity and rheometry. viscoelasticity. Applications."

# --- 4) New English "Syllabus" paragraph on row 16 --------------------------
$ws.Range("B10").Copy()
$ws.Range("B16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B16").Value = "1. Introduction. 2. Stress and deformation. 3. Types of deformation and flow of materials. 4. Fundamental equations of rheology. Flow of Newtonian and non-Newtonian fluids. 5. Viscosimetry and rheometry. 6. Rheology of dispersed systems. Colloids and emulsions. diluted solutions. Capillary viscosimetry. 7. Rheology of molten polymers. 8. Viscoelasticity. 9. Dynamic-mechanical behavior of materials. 10. Applications."

$ws.Range("C10").Copy()
$ws.Range("C16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C16").Value = "1. Introduction. 2. Stress and deformation. 3. Types of deformation and flow of materials. 4. Fundamental equations of rheology. Flow of Newtonian and non-Newtonian fluids. 5. Viscosimetry and rheometry. 6. Rheology of dispersed systems. Colloids and emulsions. diluted solutions. Capillary viscosimetry. 7. Rheology of molten polymers. 8. Viscoelasticity. 9. Dynamic-mechanical behavior of materials. 10. Applications."

Write-Host "Edits applied."
